$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume figures in place, forcing text storage
# so purely numeric-looking strings (e.g. "22.95") are not coerced
# into floating point numbers by Excel.
function Set-TextCell($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '27.416.86'
Set-TextCell 'E2' '  +0.15%  '
Set-TextCell 'D3' '1.634.77'
Set-TextCell 'E4' '  -0.03%  '
Set-TextCell 'D5' '212.31'
Set-TextCell 'E5' '  -0.72%  '
Set-TextCell 'E6' '  +4.22%  '
Set-TextCell 'E7' '  -0.05%  '
Set-TextCell 'D8' '22.95'
Set-TextCell 'E8' '  -4.33%  '
Set-TextCell 'E9' '  -2.23%  '
Set-TextCell 'D10' '0.0608'
Set-TextCell 'E10' '  -1.11%  '
Set-TextCell 'E11' '  +1.09%  '
Set-TextCell 'D12' '1.868.05'
Set-TextCell 'E12' '  -0.91%  '
Set-TextCell 'D13' '1.640.00'
Set-TextCell 'E13' '  -0.60%  '
Set-TextCell 'D14' '0.581'
Set-TextCell 'E14' '  +3.31%  '
Set-TextCell 'E15' '  -2.49%  '
Set-TextCell 'D16' '64.08'
Set-TextCell 'D17' '27.386.40'
Set-TextCell 'E17' '  +0.07%  '
Set-TextCell 'D18' '228.73'
Set-TextCell 'E18' '  -2.53%  '
Set-TextCell 'D19' '0.0₃0721'
Set-TextCell 'E19' '  -0.55%  '
Set-TextCell 'E20' '  +0.03%  '
Set-TextCell 'E21' '  +0.02%  '
Set-TextCell 'E22' '  -2.39%  '
Set-TextCell 'D23' '9.63'
Set-TextCell 'E23' '  +4.38%  '
Set-TextCell 'D24' '1.94'
Set-TextCell 'E24' '  -4.21%  '
Set-TextCell 'D25' '149.43'
Set-TextCell 'E25' '  +2.61%  '
Set-TextCell 'E26' '  -2.58%  '
Set-TextCell 'E27' '  +1.65%  '
Set-TextCell 'E28' '  -0.05%  '
Set-TextCell 'E29' '  -3.31%  '
Set-TextCell 'E30' '  -0.87%  '
Set-TextCell 'E31' '  -2.16%  '
Set-TextCell 'E32' '  -0.35%  '
Set-TextCell 'D33' '3.18'
Set-TextCell 'E33' '  +3.58%  '
Set-TextCell 'D34' '1.407.39'
Set-TextCell 'E34' '  -3.09%  '
Set-TextCell 'E35' '  +1.96%  '
Set-TextCell 'E36' '  -1.75%  '
Set-TextCell 'D38' '0.871'
Set-TextCell 'E38' '  -4.26%  '
Set-TextCell 'D39' '0.0166'
Set-TextCell 'E39' '  -1.67%  '
Set-TextCell 'E40' '  -1.17%  '
Set-TextCell 'E41' '  +0.01%  '
Set-TextCell 'D42' '0.834'
Set-TextCell 'E42' '  +5.94%  '
Set-TextCell 'E43' '  +0.72%  '
Set-TextCell 'D44' '2.22'
Set-TextCell 'E44' '  +0.35%  '
Set-TextCell 'D45' '64.65'
Set-TextCell 'E45' '  -1.35%  '
Set-TextCell 'D46' '1.776.32'
Set-TextCell 'E46' '  -0.94%  '
Set-TextCell 'D47' '1.66'
Set-TextCell 'E47' '  -3.09%  '
Set-TextCell 'D48' '85.17'
Set-TextCell 'E48' '  -3.48%  '
Set-TextCell 'E49' '  +0.24%  '
Set-TextCell 'E50' '  -1.75%  '
Set-TextCell 'D51' '7.68'
Set-TextCell 'E51' '  -1.24%  '
